# Trade #13 closed at 2026-02-16 22:58:36 - base_strategy UP +0.000%
#
# Appends a new row (row 14) to both the "All Trades" and "base_strategy"
# worksheets recording the newly-closed trade. Both sheets carry an
# identical trade log, so the same row is written to each.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 14

    $ws.Cells.Item($row, 1).Value = 13

    # The Date column ("2026-02-16") looks like a real date to Excel's
    # automatic type inference, which would silently convert it to a date
    # serial number. Force it to stay literal text (as the original export
    # stores it) with a leading apostrophe, then strip the resulting
    # quote-prefix formatting so the cell keeps the sheet's default style.
    $ws.Cells.Item($row, 2).Value = "'2026-02-16"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value = "22:58:36"
    $ws.Cells.Item($row, 4).Value = "base_strategy"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 0.5

    # Exit Price is blank (still an open-turned-closed trade with no exit
    # price recorded) - write it as a genuine empty-text cell rather than
    # leaving the cell absent, matching the source export's behaviour.
    $ws.Cells.Item($row, 7).Value = "'"
    $ws.Cells.Item($row, 7).Style = "Normal"

    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"

    # Exit Reason - also blank, same empty-text treatment as column G.
    $ws.Cells.Item($row, 16).Value = "'"
    $ws.Cells.Item($row, 16).Style = "Normal"

    $ws.Cells.Item($row, 17).Value = 0
}
